$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# In the second (manual) receptive-field table, bump in_height (C15) from 41 -> 77
$ws.Range("C15").Value = 77

# Re-enter the K/L "out_width"/"out_height" formulas for rows 15-18. Assigning the
# formula to the whole rectangular range first makes Excel group them into one
# shared formula (si="16"), matching how the author re-typed/filled the formulas.
$ws.Range("K15:L18").Formula = "=(B15-E15+2*I15)/G15+1"
$ws.Range("K15").Formula = "=(B15-E15+2*I15)/G15+1"
$ws.Range("L15").Formula = "=(C15-F15+2*J15)/H15+1"
$ws.Range("K16").Formula = "=(B16-E16+2*I16)/G16+1"
$ws.Range("L16").Formula = "=(C16-F16+2*J16)/H16+1"
$ws.Range("K17").Formula = "=(B17-E17+2*I17)/G17+1"
$ws.Range("L17").Formula = "=(C17-F17+2*J17)/H17+1"
$ws.Range("K18").Formula = "=(B18-E18+2*I18)/G18+1"
$ws.Range("L18").Formula = "=(C18-F18+2*J18)/H18+1"

# Re-enter the B/C "n_in" formulas for rows 16-18, grouped the same way (si="17").
$ws.Range("B16:C18").Formula = "=K15"
$ws.Range("B16").Formula = "=K15"
$ws.Range("C16").Formula = "=L15"
$ws.Range("B17").Formula = "=K16"
$ws.Range("C17").Formula = "=L16"
$ws.Range("B18").Formula = "=K17"
$ws.Range("C18").Formula = "=L17"

# Match the final selected cell recorded in the sheet view
$ws.Range("C15").Select()

$wb.Save()
